$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("snapshot")

$ws.Range("K2").Value = "2025-11-18T10:31:16.794522+00:00"
$ws.Range("K3").Value = "2025-11-18T10:31:16.794542+00:00"
$ws.Range("K4").Value = "2025-11-18T10:31:16.794552+00:00"
$ws.Range("K5").Value = "2025-11-18T10:31:18.775475+00:00"
$ws.Range("K6").Value = "2025-11-18T10:31:18.775509+00:00"
$ws.Range("K7").Value = "2025-11-18T10:31:18.775532+00:00"
$ws.Range("K8").Value = "2025-11-18T10:31:21.202617+00:00"
$ws.Range("K9").Value = "2025-11-18T10:31:24.025578+00:00"
$ws.Range("K10").Value = "2025-11-18T10:31:26.763751+00:00"
$ws.Range("K11").Value = "2025-11-18T10:31:26.763785+00:00"
$ws.Range("K12").Value = "2025-11-18T10:31:31.726344+00:00"
$ws.Range("K13").Value = "2025-11-18T10:31:34.028830+00:00"
$ws.Range("K14").Value = "2025-11-18T10:31:36.735706+00:00"
$ws.Range("K15").Value = "2025-11-18T10:31:36.735737+00:00"
$ws.Range("K16").Value = "2025-11-18T10:31:36.735762+00:00"
$ws.Range("K17").Value = "2025-11-18T10:31:39.002586+00:00"
$ws.Range("K18").Value = "2025-11-18T10:31:41.722696+00:00"
$ws.Range("K19").Value = "2025-11-18T10:31:44.050412+00:00"
$ws.Range("K20").Value = "2025-11-18T10:31:46.761208+00:00"
$ws.Range("K21").Value = "2025-11-18T10:31:46.761241+00:00"
$ws.Range("K22").Value = "2025-11-18T10:31:46.761261+00:00"
$ws.Range("K23").Value = "2025-11-18T10:31:46.761279+00:00"
$ws.Range("K24").Value = "2025-11-18T10:31:49.159066+00:00"
$ws.Range("K25").Value = "2025-11-18T10:31:49.159097+00:00"
$ws.Range("K26").Value = "2025-11-18T10:31:51.515759+00:00"
$ws.Range("K27").Value = "2025-11-18T10:31:51.515788+00:00"
$ws.Range("K28").Value = "2025-11-18T10:31:51.515808+00:00"
$ws.Range("K29").Value = "2025-11-18T10:31:54.808892+00:00"
$ws.Range("K30").Value = "2025-11-18T10:31:54.808923+00:00"
$ws.Range("K31").Value = "2025-11-18T10:31:57.118478+00:00"
$ws.Range("K32").Value = "2025-11-18T10:31:57.118531+00:00"
$ws.Range("K33").Value = "2025-11-18T10:31:57.118554+00:00"
$ws.Range("K34").Value = "2025-11-18T10:31:57.118575+00:00"
$ws.Range("K35").Value = "2025-11-18T10:31:57.118593+00:00"
$ws.Range("K36").Value = "2025-11-18T10:31:59.544301+00:00"
$ws.Range("K37").Value = "2025-11-18T10:31:59.544332+00:00"
$ws.Range("K38").Value = "2025-11-18T10:32:04.209999+00:00"
$ws.Range("K39").Value = "2025-11-18T10:32:04.210030+00:00"
$ws.Range("K40").Value = "2025-11-18T10:32:04.210049+00:00"
$ws.Range("K41").Value = "2025-11-18T10:32:06.458674+00:00"
